$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.365378499031067
$ws.Range("B1").Value = 2.997919797897339
$ws.Range("C1").Value = 4.882116317749023
$ws.Range("D1").Value = 1.743759393692017
$ws.Range("E1").Value = 1.103307366371155
